$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip trailing '*' (not-out marker) from HS (Highest Score) column E,
# converting these cells from text to plain numeric values.
$ws.Range("E3").Value = 100
$ws.Range("E4").Value = 2
$ws.Range("E8").Value = 1
$ws.Range("E9").Value = 62
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 9
$ws.Range("E15").Value = 58
$ws.Range("E17").Value = 76
$ws.Range("E20").Value = 8
$ws.Range("E21").Value = 100
$ws.Range("E25").Value = 1
$ws.Range("E36").Value = 39
$ws.Range("E37").Value = 59
$ws.Range("E39").Value = 68
$ws.Range("E42").Value = 84
$ws.Range("E46").Value = 64
$ws.Range("E48").Value = 25
$ws.Range("E51").Value = 17
$ws.Range("E52").Value = 48
$ws.Range("E59").Value = 59
$ws.Range("E60").Value = 37
$ws.Range("E64").Value = 8
$ws.Range("E66").Value = 10
$ws.Range("E67").Value = 63
$ws.Range("E71").Value = 95
$ws.Range("E72").Value = 13
$ws.Range("E73").Value = 84
$ws.Range("E77").Value = 33
$ws.Range("E78").Value = 75
$ws.Range("E79").Value = 54
$ws.Range("E80").Value = 10
$ws.Range("E81").Value = 20
$ws.Range("E82").Value = 54
$ws.Range("E86").Value = 107
$ws.Range("E87").Value = 91
$ws.Range("E89").Value = 124
$ws.Range("E91").Value = 15
$ws.Range("E98").Value = 12
$ws.Range("E100").Value = 66
$ws.Range("E102").Value = 63
$ws.Range("E104").Value = 25
$ws.Range("E105").Value = 47
$ws.Range("E107").Value = 36
$ws.Range("E108").Value = 88
$ws.Range("E110").Value = 66
$ws.Range("E113").Value = 62
$ws.Range("E115").Value = 106

# Update the saved selection/active cell
$ws.Range("E3").Select() | Out-Null
